$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "CAWI"
$ws.Range("B15").Value = "Computer Assissted Web Interviewing"

$ws.Range("F14").Select()
